# feat: add terms billing model property to product item
#
# Items sheet: the old single "Billing Frequency" column is split into two
# columns - a new "Billing Model" column (quantity/usage) inserted right
# before the existing frequency column, which itself is relabeled
# "Billing Period".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# Insert a new column before the old "Billing Frequency" column (G),
# shifting it (and everything after it) one column to the right.
$ws.Columns.Item(7).Insert()

# Give the new column roughly the same width as the Description column,
# matching the wide free-text style used by its neighbour.
$ws.Columns.Item(7).ColumnWidth = 68.83

# Header row: new column G is "Billing Model"; the old "Billing Frequency"
# header, now shifted into column H, becomes "Billing Period".
$ws.Cells.Item(1, 7).Value = "Billing Model"
$ws.Cells.Item(1, 8).Value = "Billing Period"

# Populate the Billing Model values for each item row.
$ws.Cells.Item(2, 7).Value = "quantity"
$ws.Cells.Item(3, 7).Value = "usage"
$ws.Cells.Item(4, 7).Value = "usage"
$ws.Cells.Item(5, 7).Value = "usage"
$ws.Cells.Item(6, 7).Value = "usage"
$ws.Cells.Item(7, 7).Value = "usage"

# Reflect the saved selection/scroll state from the author's session.
$ws.Range("H11").Select()
